# Apply updated crypto price/volume data (and the Stellar/ICP row swap)
# to match the latest scraped values from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: some "Price" values look like plain numbers (e.g. 250.11).
# Prefixing with a single quote forces Excel to store them as text
# (matching the inlineStr/text cells in the workbook) instead of
# converting them to floating point numbers.

$ws.Range("D2").Value = '41.852.62'
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("D3").Value = '2.228.21'
$ws.Range("E3").Value = '  -0.86%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''250.11'
$ws.Range("E5").Value = '  +7.05%  '
$ws.Range("D6").Value = '''0.629'
$ws.Range("E6").Value = '  -0.97%  '
$ws.Range("D7").Value = '''71.51'
$ws.Range("E7").Value = '  +2.15%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").Value = '''0.585'
$ws.Range("E9").Value = '  +3.93%  '
$ws.Range("D10").Value = '''41.36'
$ws.Range("E10").Value = '  +14.86%  '
$ws.Range("E11").Value = '  -2.62%  '
$ws.Range("D12").Value = '''58.29'
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("E13").Value = '  +1.22%  '
$ws.Range("D14").Value = '''7.01'
$ws.Range("E14").Value = '  +2.57%  '
$ws.Range("D15").Value = '2.560.36'
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("D16").Value = '''14.97'
$ws.Range("E16").Value = '  -0.94%  '
$ws.Range("D17").Value = '''0.864'
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("D18").Value = '2.229.49'
$ws.Range("E18").Value = '  -0.88%  '
$ws.Range("D19").Value = '41.758.83'
$ws.Range("E19").Value = '  -0.88%  '
$ws.Range("D20").Value = '0.0₃0972'
$ws.Range("E20").Value = '  -0.68%  '
$ws.Range("D21").Value = '''6.21'
$ws.Range("E21").Value = '  -0.88%  '
$ws.Range("D22").Value = '''72.80'
$ws.Range("E22").Value = '  -0.95%  '
$ws.Range("D23").Value = '''235.36'
$ws.Range("E23").Value = '  -0.85%  '
$ws.Range("D24").Value = '''2.16'
$ws.Range("E24").Value = '  +6.50%  '
$ws.Range("D25").Value = '''4.21'
$ws.Range("E25").Value = '  +15.21%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").Value = '''2.53'
$ws.Range("E27").Value = '  +6.94%  '
$ws.Range("D28").Value = '''10.55'
$ws.Range("E28").Value = '  +4.54%  '
$ws.Range("E29").Value = '  +0.82%  '
$ws.Range("D30").Value = '''171.39'
$ws.Range("E30").Value = '  +1.21%  '
$ws.Range("D31").Value = '''20.74'
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("D32").Value = '''0.122'
$ws.Range("E32").Value = '  +2.10%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = '''5.65'
$ws.Range("E33").Value = '  +4.77%  '
$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D34").Value = '''0.125'
$ws.Range("E34").Value = '  -2.15%  '
$ws.Range("E35").Value = '  +0.71%  '
$ws.Range("E36").Value = '  -1.91%  '
$ws.Range("D37").Value = '''26.13'
$ws.Range("E37").Value = '  +17.41%  '
$ws.Range("D38").Value = '''3.93'
$ws.Range("E38").Value = '  +8.08%  '
$ws.Range("D39").Value = '''0.0304'
$ws.Range("E39").Value = '  +13.81%  '
$ws.Range("E40").Value = '  +1.16%  '
$ws.Range("D41").Value = '''68.38'
$ws.Range("E41").Value = '  +2.42%  '
$ws.Range("E42").Value = '  -1.65%  '
$ws.Range("D43").Value = '''11.85'
$ws.Range("E43").Value = '  +16.59%  '
$ws.Range("D44").Value = '''0.208'
$ws.Range("E44").Value = '  +7.53%  '
$ws.Range("D45").Value = '''4.88'
$ws.Range("E45").Value = '  -2.64%  '
$ws.Range("D46").Value = '''8.83'
$ws.Range("E46").Value = '  -2.28%  '
$ws.Range("D47").Value = '''4.76'
$ws.Range("E47").Value = '  +8.39%  '
$ws.Range("E48").Value = '  +0.98%  '
$ws.Range("E49").Value = '  +0.25%  '
$ws.Range("E50").Value = '  +7.25%  '
$ws.Range("E51").Value = '  +0.61%  '
